$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 271.26086
$ws.Range("J2").Value = 578
$ws.Range("L2").Value = 578
$ws.Range("N2").Value = -804
$ws.Range("H17").Value = 1193.963
$ws.Range("J17").Value = 1209.48
$ws.Range("L17").Value = 3628.44
$ws.Range("N17").Value = -3964.44
$ws.Range("H43").Value = 999
$ws.Range("J43").Value = 999
$ws.Range("L43").Value = 999
$ws.Range("N43").Value = -1137
$ws.Range("H70").Value = 71402.7
$ws.Range("J70").Value = 75809
$ws.Range("L70").Value = 227427
$ws.Range("N70").Value = -227967
$ws.Range("H73").Value = 71402.7
$ws.Range("J73").Value = 75809
$ws.Range("L73").Value = 227427
$ws.Range("N73").Value = -229299
$ws.Range("H100").Value = 1900.4286
$ws.Range("I100").Value = 474.5
$ws.Range("J100").Value = 2470.8
$ws.Range("K100").Value = 474.5
$ws.Range("L100").Value = 2470.8
$ws.Range("M100").Value = 66.5
$ws.Range("N100").Value = -3552.8
$ws.Range("H107").Value = 364.25
$ws.Range("I107").Value = 364.25
$ws.Range("K107").Value = 364.25
$ws.Range("M107").Value = 1555.75
$ws.Range("H112").Value = 2808.875
$ws.Range("J112").Value = 3128.5
$ws.Range("L112").Value = 9385.5
$ws.Range("N112").Value = -11601.5
$ws.Range("H113").Value = 4014.1667
$ws.Range("I113").Value = 3917
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 3917
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -663
$ws.Range("N113").Value = -11008

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1296.6666
$ws.Range("I2").Value = 695
$ws.Range("K2").Value = 695
$ws.Range("M2").Value = -582
$ws.Range("H43").Value = 34247.5
$ws.Range("J43").Value = 31797.6
$ws.Range("L43").Value = 31797.6
$ws.Range("N43").Value = -32423.6
$ws.Range("H74").Value = 1676.8182
$ws.Range("I74").Value = 1684.5
$ws.Range("K74").Value = 1684.5
$ws.Range("M74").Value = -810.5
$ws.Range("H77").Value = 1676.8182
$ws.Range("I77").Value = 1684.5
$ws.Range("K77").Value = 8422.5
$ws.Range("M77").Value = -4054.5
$ws.Range("H110").Value = 12334666
$ws.Range("I110").Value = 12334666
$ws.Range("K110").Value = 12334666
$ws.Range("M110").Value = -12332621
$ws.Range("H116").Value = 1296.6666
$ws.Range("I116").Value = 695
$ws.Range("K116").Value = 695
$ws.Range("M116").Value = 1599

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1296.6666
$ws.Range("I3").Value = 695
$ws.Range("K3").Value = 695
$ws.Range("M3").Value = -581

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H99").Value = 7086.6
$ws.Range("I99").Value = 1358.5
$ws.Range("J99").Value = 29999
$ws.Range("K99").Value = 1358.5
$ws.Range("L99").Value = 29999
$ws.Range("M99").Value = 139.5
$ws.Range("N99").Value = -32995
$ws.Range("H126").Value = 7086.6
$ws.Range("I126").Value = 1358.5
$ws.Range("J126").Value = 29999
$ws.Range("K126").Value = 4075.5
$ws.Range("L126").Value = 89997
$ws.Range("M126").Value = -1605.5
$ws.Range("N126").Value = -94937
$ws.Range("H133").Value = 34151.625
$ws.Range("I133").Value = 25296
$ws.Range("K133").Value = 25296
$ws.Range("M133").Value = -22766

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 392.81818
$ws.Range("I7").Value = 369.2
$ws.Range("J7").Value = 412.5
$ws.Range("K7").Value = 1107.6
$ws.Range("L7").Value = 1237.5
$ws.Range("M7").Value = -995.5999999999999
$ws.Range("N7").Value = -1461.5
$ws.Range("H12").Value = 1552.44
$ws.Range("I12").Value = 930.5454999999999
$ws.Range("J12").Value = 2041.0714
$ws.Range("K12").Value = 2791.6365
$ws.Range("L12").Value = 6123.2142
$ws.Range("M12").Value = -2618.6365
$ws.Range("N12").Value = -6469.2142
$ws.Range("H46").Value = 1067.9
$ws.Range("I46").Value = 1048.625
$ws.Range("J46").Value = 1145
$ws.Range("K46").Value = 3145.875
$ws.Range("L46").Value = 3435
$ws.Range("M46").Value = -3054.875
$ws.Range("N46").Value = -3617
$ws.Range("H56").Value = 14471.583
$ws.Range("I56").Value = 14471.583
$ws.Range("K56").Value = 14471.583
$ws.Range("M56").Value = -13941.583
$ws.Range("H92").Value = 281.4
$ws.Range("I92").Value = 247.625
$ws.Range("K92").Value = 742.875
$ws.Range("M92").Value = 505.125
$ws.Range("H128").Value = 647557.2
$ws.Range("I128").Value = 647557.2
$ws.Range("K128").Value = 1942671.6
$ws.Range("M128").Value = -1937691.6
$ws.Range("H133").Value = 2299.5
$ws.Range("I133").Value = 2299.5
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 6898.5
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -1838.5
$ws.Range("N133").ClearContents()
$ws.Range("H140").Value = 11735.066
$ws.Range("I140").Value = 1333
$ws.Range("J140").Value = 12478.071
$ws.Range("K140").Value = 3999
$ws.Range("L140").Value = 37434.213
$ws.Range("M140").Value = 1181
$ws.Range("N140").Value = -47794.213

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3552.261
$ws.Range("I122").Value = 3485.15
$ws.Range("J122").Value = 3999.6667
$ws.Range("K122").Value = 10455.45
$ws.Range("L122").Value = 11999.0001
$ws.Range("M122").Value = -8005.450000000001
$ws.Range("N122").Value = -16899.0001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 58663.668
$ws.Range("I46").Value = 101994.6
$ws.Range("J46").Value = 4500
$ws.Range("K46").Value = 101994.6
$ws.Range("L46").Value = 4500
$ws.Range("M46").Value = -101806.6
$ws.Range("N46").Value = -4876
$ws.Range("H55").Value = 234.75
$ws.Range("I55").Value = 210.25
$ws.Range("J55").Value = 259.25
$ws.Range("K55").Value = 210.25
$ws.Range("L55").Value = 259.25
$ws.Range("M55").Value = -37.25
$ws.Range("N55").Value = -605.25
$ws.Range("H74").Value = 25065.666
$ws.Range("I74").Value = 24598.5
$ws.Range("J74").Value = 26000
$ws.Range("K74").Value = 24598.5
$ws.Range("L74").Value = 26000
$ws.Range("M74").Value = -23600.5
$ws.Range("N74").Value = -27996
$ws.Range("H77").Value = 25065.666
$ws.Range("I77").Value = 24598.5
$ws.Range("J77").Value = 26000
$ws.Range("K77").Value = 73795.5
$ws.Range("L77").Value = 78000
$ws.Range("M77").Value = -68803.5
$ws.Range("N77").Value = -87984
$ws.Range("H122").Value = 6180.6294
$ws.Range("I122").Value = 6030.9546
$ws.Range("J122").Value = 6839.2
$ws.Range("K122").Value = 18092.8638
$ws.Range("L122").Value = 20517.6
$ws.Range("M122").Value = -15642.8638
$ws.Range("N122").Value = -25417.6
$ws.Range("H139").Value = 27500
$ws.Range("I139").Value = 27500
$ws.Range("K139").Value = 27500
$ws.Range("M139").Value = -22360

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 21062.5
$ws.Range("I62").Value = 23751
$ws.Range("J62").Value = 20166.334
$ws.Range("K62").Value = 23751
$ws.Range("L62").Value = 20166.334
$ws.Range("M62").Value = -23127
$ws.Range("N62").Value = -21414.334
$ws.Range("H65").Value = 21062.5
$ws.Range("I65").Value = 23751
$ws.Range("J65").Value = 20166.334
$ws.Range("K65").Value = 118755
$ws.Range("L65").Value = 100831.67
$ws.Range("M65").Value = -115635
$ws.Range("N65").Value = -107071.67
$ws.Range("H81").Value = 715634.3
$ws.Range("I81").Value = 1452.1538
$ws.Range("J81").Value = 10000002
$ws.Range("K81").Value = 2904.3076
$ws.Range("L81").Value = 20000004
$ws.Range("M81").Value = -1843.3076
$ws.Range("N81").Value = -20002126
$ws.Range("H84").Value = 715634.3
$ws.Range("I84").Value = 1452.1538
$ws.Range("J84").Value = 10000002
$ws.Range("K84").Value = 14521.538
$ws.Range("L84").Value = 100000020
$ws.Range("M84").Value = -9217.538
$ws.Range("N84").Value = -100010628
$ws.Range("H100").Value = 5883323
$ws.Range("I100").Value = 7693314
$ws.Range("J100").Value = 852
$ws.Range("K100").Value = 15386628
$ws.Range("L100").Value = 1704
$ws.Range("M100").Value = -15386087
$ws.Range("N100").Value = -2786
$ws.Range("H119").Value = 120000
$ws.Range("J119").Value = 120000
$ws.Range("L119").Value = 120000
$ws.Range("N119").Value = -129676

Write-Host "Applied all Sagittarius_Profits edits"